$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" footer field text from
#    9/5/2017 -> 9/10/2018 on the slide master and every custom (slide)
#    layout that carries a Date Placeholder.
# ---------------------------------------------------------------------------
$oldDate = "9/5/2017"
$newDate = "9/10/2018"

$sm = $p.SlideMaster

for ($i = 1; $i -le $sm.Shapes.Count; $i++) {
    $sh = $sm.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$layouts = $sm.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 1, "Content Placeholder 2": the two runs
#       "     " + "(can be done on calculator): "
#    become a single run "     (can be done on calculator): ".
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$content = $slide1.Shapes.Item(2)
$tr = $content.TextFrame.TextRange
$mergedText = "     (can be done on calculator): "

for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    $paraText = $para.Text.TrimEnd([char]13)
    if ($paraText -eq $mergedText) {
        $len = $para.Length
        if ($para.Text.EndsWith([string][char]13)) {
            $len = $len - 1
        }
        $whole = $tr.Characters($para.Start, $len)
        $whole.Text = $mergedText
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Slide 1: remove the "TextBox 2" shape that reads
#    "I promise we won't spend the entire time on this again!"
# ---------------------------------------------------------------------------
for ($i = $slide1.Shapes.Count; $i -ge 1; $i--) {
    $sh = $slide1.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 2") {
        $sh.Delete()
    }
}
